$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row value updates (Price / Volume(1h) columns) ---
$ws.Range("D2").Value = '29.658.01'
$ws.Range("E2").Value = '  +2.78%  '

$ws.Range("D3").Value = '1.864.64'
$ws.Range("E3").Value = '  +2.08%  '

$ws.Range("D5").Value = '246.24'
$ws.Range("E5").Value = '  +2.82%  '

$ws.Range("D6").Value = '0.7008'
$ws.Range("E6").Value = '  +1.98%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '0.07766'
$ws.Range("E8").Value = '  +1.95%  '

$ws.Range("D9").Value = '0.3080'
$ws.Range("E9").Value = '  +2.03%  '

$ws.Range("D10").Value = '23.71'
$ws.Range("E10").Value = '  +1.11%  '

$ws.Range("D11").Value = '0.07821'
$ws.Range("E11").Value = '  +1.19%  '

$ws.Range("D14").Value = '92.91'
$ws.Range("E14").Value = '  +3.19%  '

$ws.Range("D15").Value = '0.6956'
$ws.Range("E15").Value = '  +3.61%  '

$ws.Range("D16").Value = '6.627'
$ws.Range("E16").Value = '  +3.42%  '

$ws.Range("D17").Value = '29.666.33'
$ws.Range("E17").Value = '  +2.78%  '

$ws.Range("D18").Value = '0.000008382'
$ws.Range("E18").Value = '  +1.41%  '

$ws.Range("D19").Value = '2.112.61'
$ws.Range("E19").Value = '  +1.57%  '

$ws.Range("D20").Value = '243.66'
$ws.Range("E20").Value = '  +0.52%  '

$ws.Range("D21").Value = '12.81'
$ws.Range("E21").Value = '  +1.68%  '

$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("D23").Value = '7.647'
$ws.Range("E23").Value = '  +3.45%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").Value = '0.1526'
$ws.Range("E25").Value = '  +4.01%  '

$ws.Range("D26").Value = '8.961'
$ws.Range("E26").Value = '  +3.06%  '

$ws.Range("E27").Value = '  +0.00%  '

$ws.Range("E28").Value = '  +1.48%  '

$ws.Range("D29").Value = '1.549'
$ws.Range("E29").Value = '  +1.32%  '

$ws.Range("D30").Value = '4.276'
$ws.Range("E30").Value = '  +1.93%  '

$ws.Range("D31").Value = '4.209'
$ws.Range("E31").Value = '  +1.61%  '

$ws.Range("D32").Value = '1.199'
$ws.Range("E32").Value = '  +0.72%  '

$ws.Range("E33").Value = '  +0.45%  '

$ws.Range("D34").Value = '0.7884'
$ws.Range("E34").Value = '  +4.18%  '

$ws.Range("D35").Value = '1.913'
$ws.Range("E35").Value = '  +5.72%  '

$ws.Range("D36").Value = '1.159'
$ws.Range("E36").Value = '  +1.86%  '

$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("D38").Value = '1.343.46'
$ws.Range("E38").Value = '  +10.98%  '

$ws.Range("D39").Value = '0.01884'
$ws.Range("E39").Value = '  +2.99%  '

$ws.Range("D40").Value = '2.743'
$ws.Range("E40").Value = '  +2.57%  '

$ws.Range("D41").Value = '0.9596'
$ws.Range("E41").Value = '  +5.08%  '

$ws.Range("D42").Value = '6.029'
$ws.Range("E42").Value = '  +14.21%  '

$ws.Range("D43").Value = '106.77'
$ws.Range("E43").Value = '  -1.59%  '

$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("D47").Value = '2.012.59'
$ws.Range("E47").Value = '  +1.27%  '

$ws.Range("D48").Value = '65.25'
$ws.Range("E48").Value = '  +4.52%  '

$ws.Range("D49").Value = '0.5216'
$ws.Range("E49").Value = '  +1.06%  '

$ws.Range("D50").Value = '1.791'
$ws.Range("E50").Value = '  +4.10%  '

$ws.Range("D51").Value = '7.025'
$ws.Range("E51").Value = '  +2.62%  '

# --- Row 12 / Row 13 swap: Polkadot and WrappedEther switch ranking ---
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "5.174"
$ws.Range("E12").Value = "  +2.66%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.866.17"
$ws.Range("E13").Value = "  +2.15%  "

# --- Row 45 / Row 46 swap: EnergySwap and BabyDogeCoin switch ranking ---
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "9.805"
$ws.Range("E45").Value = "  +4.12%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000125"
$ws.Range("E46").Value = "  +2.23%  "
